# The workbook has a single worksheet "quadratic-svm-score" with a small
# prediction table:
#   Row | 1-s__Clostridium_AP scindens | prediction
#   even_MAG-GUT11004.fa | 1 | 1
#   even_MAG-GUT14745.fa | 1 | 1
#
# The "single child" (even_MAG-GUT14745.fa) output was removed from this
# run, and the remaining child's prediction score was updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quadratic-svm-score")

# Update the prediction score for the remaining row (even_MAG-GUT11004.fa)
$ws.Range("B2").Value = 139.0436161570187

# Drop the even_MAG-GUT14745.fa row entirely - it is no longer an output
$ws.Rows.Item(3).Delete()

# Keep the surviving header/label cells as text-formatted, same as before
$ws.Range("A1:C1").NumberFormat = "@"
$ws.Range("A2").NumberFormat = "@"
